# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages have now been received:
#   - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview sheet and on each language sheet.
#   - Each language sheet's "Latest Target File" / "Latest Handback File" columns
#     (I, J) get populated with the handed-back file names (with hyperlinks, like
#     column A), and the de-de sheet also receives a "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdName1 = "9c2f69b4-dffe-42af-8777-0775aede5eda.md"
$mdName2 = "f4c95d07-63ab-4881-a11c-d773bcb16a0a.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40d45debf186365b6f15006088a71cb9e206f6a3/e2e/9c2f69b4-dffe-42af-8777-0775aede5eda.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40d45debf186365b6f15006088a71cb9e206f6a3/e2e/f4c95d07-63ab-4881-a11c-d773bcb16a0a.md"

$zhXlf1 = "9c2f69b4-dffe-42af-8777-0775aede5eda.be54b896dd9c4f25f9cce665c9b3690ef046edd4.zh-cn.xlf"
$zhXlf2 = "f4c95d07-63ab-4881-a11c-d773bcb16a0a.78fec0aa71e06d7cebdb642d5bb427280b7e1391.zh-cn.xlf"
$deXlf1 = "9c2f69b4-dffe-42af-8777-0775aede5eda.be54b896dd9c4f25f9cce665c9b3690ef046edd4.de-de.xlf"
$deXlf2 = "f4c95d07-63ab-4881-a11c-d773bcb16a0a.78fec0aa71e06d7cebdb642d5bb427280b7e1391.de-de.xlf"

# Excel's ColumnWidth setter snaps to a 1/6-character pixel grid internally
# (stored_width = round(input*6)/6 + 5/6), so request the value that lands
# closest to the desired rendered width.
function Set-ColWidth($ws, $colIndex, $target) {
    $n = [Math]::Round(($target - (5.0/6.0)) * 6.0)
    $input = $n / 6.0
    $ws.Columns.Item($colIndex).ColumnWidth = $input
}

# ---------------------------------------------------------------------------
# Overview sheet: the Status columns (zh-cn = E, de-de = F) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both rows.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

Set-ColWidth $overview 5 29.9777047293527
Set-ColWidth $overview 6 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("I2").Value = $mdName1
$zh.Range("J2").Value = $zhXlf1
$zh.Range("K2").Value = "2016-09-01 12:32:33"

$zh.Range("I3").Value = $mdName2
$zh.Range("J3").Value = $zhXlf2
$zh.Range("K3").Value = "2016-09-01 12:32:33"

# Rebuild hyperlinks so that the new "Latest Target File" links (I2, I3) sit
# alongside the existing "Source File Name" links (A2, A3).
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $mdUrl1, "", "", $mdName1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), $mdUrl2, "", "", $mdName2) | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null

Set-ColWidth $zh 3 29.9777047293527
Set-ColWidth $zh 9 40
Set-ColWidth $zh 10 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("I2").Value = $mdName1
$de.Range("J2").Value = $deXlf1
$de.Range("K2").Value = "2016-09-01 12:32:40"

$de.Range("I3").Value = $mdName2
$de.Range("J3").Value = $deXlf2
$de.Range("K3").Value = "2016-09-01 12:32:40"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $mdUrl1, "", "", $mdName1) | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), $mdUrl2, "", "", $mdName2) | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null

Set-ColWidth $de 3 29.9777047293527
Set-ColWidth $de 9 40
Set-ColWidth $de 10 40

Write-Host "Handback report generated."
